$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header
# formatting (bold / centered / bordered) used by the rest of row 1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new columns I (I0) and J (IF), rows 2-15.
$values = @(
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(5, 6),
    @(9, 9),
    @(6, 6),
    @(4, 4),
    @(9, 9),
    @(8, 8),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
